$wb = $excel.ActiveWorkbook

# Sheet "high_loadings" - update Category (column B) values
$ws1 = $wb.Worksheets.Item("high_loadings")
$ws1.Cells.Item(2, 2).Value = 2
$ws1.Cells.Item(3, 2).Value = 2
$ws1.Cells.Item(4, 2).Value = 3
$ws1.Cells.Item(5, 2).Value = 2
$ws1.Cells.Item(6, 2).Value = 3
$ws1.Cells.Item(7, 2).Value = 1
$ws1.Cells.Item(8, 2).Value = 3
$ws1.Cells.Item(9, 2).Value = 1
$ws1.Cells.Item(11, 2).Value = 3
$ws1.Cells.Item(14, 2).Value = 2
$ws1.Cells.Item(15, 2).Value = 3
$ws1.Cells.Item(17, 2).Value = 1
$ws1.Cells.Item(24, 2).Value = 2
$ws1.Cells.Item(27, 2).Value = 1
$ws1.Cells.Item(30, 2).Value = 3
$ws1.Cells.Item(32, 2).Value = 3
$ws1.Cells.Item(34, 2).Value = 2
$ws1.Cells.Item(35, 2).Value = 1
$ws1.Cells.Item(36, 2).Value = 3
$ws1.Cells.Item(37, 2).Value = 2
$ws1.Cells.Item(40, 2).Value = 1
$ws1.Cells.Item(42, 2).Value = 2
$ws1.Cells.Item(44, 2).Value = 1
$ws1.Cells.Item(47, 2).Value = 2
$ws1.Cells.Item(50, 2).Value = 2
$ws1.Cells.Item(51, 2).Value = 3
$ws1.Cells.Item(52, 2).Value = 2
$ws1.Cells.Item(55, 2).Value = 2
$ws1.Cells.Item(56, 2).Value = 2
$ws1.Cells.Item(58, 2).Value = 1
$ws1.Cells.Item(59, 2).Value = 2
$ws1.Cells.Item(60, 2).Value = 1
$ws1.Cells.Item(61, 2).Value = 1
$ws1.Cells.Item(63, 2).Value = 2
$ws1.Cells.Item(64, 2).Value = 3
$ws1.Cells.Item(65, 2).Value = 2
$ws1.Cells.Item(66, 2).Value = 3
$ws1.Cells.Item(67, 2).Value = 1
$ws1.Cells.Item(70, 2).Value = 3

# Sheet "Influencers_uniques" - update Category (column B) values
$ws2 = $wb.Worksheets.Item("Influencers_uniques")
$ws2.Cells.Item(2, 2).Value = 3
$ws2.Cells.Item(4, 2).Value = 1
$ws2.Cells.Item(5, 2).Value = 2
$ws2.Cells.Item(6, 2).Value = 2
$ws2.Cells.Item(7, 2).Value = 2
$ws2.Cells.Item(8, 2).Value = 1
$ws2.Cells.Item(9, 2).Value = 3
$ws2.Cells.Item(10, 2).Value = 2
$ws2.Cells.Item(12, 2).Value = 1
$ws2.Cells.Item(15, 2).Value = 3
$ws2.Cells.Item(16, 2).Value = 1
$ws2.Cells.Item(17, 2).Value = 2
$ws2.Cells.Item(21, 2).Value = 1
$ws2.Cells.Item(25, 2).Value = 2
$ws2.Cells.Item(27, 2).Value = 3
$ws2.Cells.Item(28, 2).Value = 1
$ws2.Cells.Item(31, 2).Value = 2
$ws2.Cells.Item(32, 2).Value = 2
$ws2.Cells.Item(34, 2).Value = 1
$ws2.Cells.Item(38, 2).Value = 3
$ws2.Cells.Item(39, 2).Value = 2
$ws2.Cells.Item(40, 2).Value = 1
$ws2.Cells.Item(41, 2).Value = 2
$ws2.Cells.Item(42, 2).Value = 1
$ws2.Cells.Item(45, 2).Value = 3
$ws2.Cells.Item(47, 2).Value = 2
$ws2.Cells.Item(48, 2).Value = 1
$ws2.Cells.Item(49, 2).Value = 2
$ws2.Cells.Item(50, 2).Value = 2
$ws2.Cells.Item(51, 2).Value = 2
$ws2.Cells.Item(56, 2).Value = 3
$ws2.Cells.Item(59, 2).Value = 1
$ws2.Cells.Item(60, 2).Value = 3
$ws2.Cells.Item(62, 2).Value = 1
$ws2.Cells.Item(63, 2).Value = 2
$ws2.Cells.Item(64, 2).Value = 1
$ws2.Cells.Item(65, 2).Value = 3
$ws2.Cells.Item(66, 2).Value = 3
$ws2.Cells.Item(68, 2).Value = 2
$ws2.Cells.Item(70, 2).Value = 3
